# Aanvulling lijst To Do
# Fill in row 14 of the ToDO table with a new finished task ("Aanpassing CSS
# naar afspraken") and move the active selection to A15, the next empty row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ToDO")

$ws.Range("A14").Value = "Aanpassing CSS naar afspraken"
$ws.Range("B14").Value = "15 minuten"
$ws.Range("C14").Value = "30 minuten"

# D14 is a brand-new cell (row 14 previously had no D cell at all). Copy the
# existing short-date format from the cell right above it so the new cell
# reuses the same style record (numFmtId 14) instead of minting a new one,
# then fill in the date serial value itself.
$ws.Range("D13").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = 41365

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "Steven V"
$ws.Range("G14").Value = "Solved"
$ws.Range("H14").Value = "CSS"

$ws.Range("A15").Select()
